$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header values in P1 and Q1, copying the formatting (bold font,
# centered/top aligned, thin box border) from the adjacent header cell O1.
$ws.Range("O1").Copy()
$ws.Range("P1:Q1").PasteSpecial(-4122)
$ws.Range("P1:Q1").Value = 14
$ws.Range("Q1").Value = 15

for ($r = 2; $r -le 25; $r++) {
    $ws.Cells.Item($r, 9).Value = 2   # I: 1 -> 2
    $ws.Cells.Item($r, 11).Value = 1  # K: 2 -> 1
    $ws.Cells.Item($r, 13).Value = 2  # M: 1 -> 2
    $ws.Cells.Item($r, 15).Value = 1  # O: 2 -> 1
    $ws.Cells.Item($r, 16).Value = 2  # P: new column = 2
    $ws.Cells.Item($r, 17).Value = 2  # Q: new column = 2
}
